$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.7435299789381
$ws.Range("D2").Value = 7.895511040167776
$ws.Range("E2").Value = 23.97727777931222
$ws.Range("F2").Value = 45.93636939641061
$ws.Range("G2").Value = 3.635167101352711
$ws.Range("I2").Value = 24.60619259076782
$ws.Range("L2").Value = 13.36631457433549
$ws.Range("N2").Value = 17.96673420643691

$ws.Range("B3").Value = 17.36251559319772
$ws.Range("D3").Value = 7.777055799274089
$ws.Range("E3").Value = 22.46871719913002
$ws.Range("F3").Value = 44.56221823164615
$ws.Range("G3").Value = 3.643944028944533
$ws.Range("I3").Value = 24.7855648676278
$ws.Range("L3").Value = 13.07496835463885
$ws.Range("N3").Value = 18.04296343812682

$ws.Range("B4").Value = 17.12891609210078
$ws.Range("D4").Value = 7.707388924229538
$ws.Range("E4").Value = 21.48973561322724
$ws.Range("F4").Value = 43.7203171692588
$ws.Range("G4").Value = 3.649583682778916
$ws.Range("I4").Value = 24.90233406872711
$ws.Range("L4").Value = 12.89739030964981
$ws.Range("N4").Value = 18.09191414805754

$ws.Range("B5").Value = 17.03394399529442
$ws.Range("D5").Value = 7.679797929502742
$ws.Range("E5").Value = 21.07755138501533
$ws.Range("F5").Value = 43.37826900217194
$ws.Range("G5").Value = 3.651945353778197
$ws.Range("I5").Value = 24.95158007281294
$ws.Range("L5").Value = 12.82545700682666
$ws.Range("N5").Value = 18.11240481898043

$ws.Range("B6").Value = 17.01819132923751
$ws.Range("D6").Value = 7.675265485536334
$ws.Range("E6").Value = 21.00830955298322
$ws.Range("F6").Value = 43.3215517292835
$ws.Range("G6").Value = 3.652341353881346
$ws.Range("I6").Value = 24.95985750790918
$ws.Range("L6").Value = 12.81354157721302
$ws.Range("N6").Value = 18.11584016363455

$ws.Range("B7").Value = 17.127634187703
$ws.Range("D7").Value = 7.707013552693259
$ws.Range("E7").Value = 21.48423033371274
$ws.Range("F7").Value = 43.71569923138499
$ws.Range("G7").Value = 3.649615275526389
$ws.Range("I7").Value = 24.90299149805998
$ws.Range("L7").Value = 12.89641831066488
$ws.Range("N7").Value = 18.09218828998938

$ws.Range("B8").Value = 17.61216461540456
$ws.Range("D8").Value = 7.854049054748264
$ws.Range("E8").Value = 23.46807210705305
$ws.Range("F8").Value = 45.46251755401737
$ws.Range("G8").Value = 3.638141652905393
$ws.Range("I8").Value = 24.6666592409657
$ws.Range("L8").Value = 13.26564723809274
$ws.Range("N8").Value = 17.99257472080827

$ws.Range("B9").Value = 18.55931637135392
$ws.Range("D9").Value = 8.165368044127973
$ws.Range("E9").Value = 26.94021536914255
$ws.Range("F9").Value = 48.878773663556
$ws.Range("G9").Value = 3.617609265706515
$ws.Range("I9").Value = 24.25615855641836
$ws.Range("L9").Value = 13.99571565968178
$ws.Range("N9").Value = 17.81410569474417

$ws.Range("B10").Value = 19.24599174251404
$ws.Range("D10").Value = 8.406205939829963
$ws.Range("E10").Value = 29.23739386317487
$ws.Range("F10").Value = 51.35350858330505
$ws.Range("G10").Value = 3.603693866304828
$ws.Range("I10").Value = 23.98729238002272
$ws.Range("L10").Value = 14.53010945152119
$ws.Range("N10").Value = 17.69305686601538

$ws.Range("B11").Value = 19.55493969138293
$ws.Range("D11").Value = 8.517976642005531
$ws.Range("E11").Value = 30.22775824756393
$ws.Range("F11").Value = 52.46617914502883
$ws.Range("G11").Value = 3.597610879038639
$ws.Range("I11").Value = 23.87219334128898
$ws.Range("L11").Value = 14.77167187840343
$ws.Range("N11").Value = 17.64012931143575

$ws.Range("B12").Value = 19.6713268742187
$ws.Range("D12").Value = 8.560583658638789
$ws.Range("E12").Value = 30.5949688065304
$ws.Range("F12").Value = 52.88522932117169
$ws.Range("G12").Value = 3.595342429210892
$ws.Range("I12").Value = 23.82965519622438
$ws.Range("L12").Value = 14.86283813174334
$ws.Range("N12").Value = 17.62039080471944

$ws.Range("B13").Value = 19.64628941554291
$ws.Range("D13").Value = 8.551395511864762
$ws.Range("E13").Value = 30.51623054043891
$ws.Range("F13").Value = 52.79508749988742
$ws.Range("G13").Value = 3.595829429897199
$ws.Range("I13").Value = 23.83876978836221
$ws.Range("L13").Value = 14.84321888517405
$ws.Range("N13").Value = 17.62462837980736

$ws.Range("B14").Value = 19.5645276710638
$ws.Range("D14").Value = 8.52147650396703
$ws.Range("E14").Value = 30.25812535124355
$ws.Range("F14").Value = 52.50070283779761
$ws.Range("G14").Value = 3.597423552851442
$ws.Range("I14").Value = 23.86867264846518
$ws.Range("L14").Value = 14.7791788377419
$ws.Range("N14").Value = 17.63849933974938

$ws.Range("B15").Value = 19.51436428489735
$ws.Range("D15").Value = 8.503185886592373
$ws.Range("E15").Value = 30.09901141405073
$ws.Range("F15").Value = 52.32007349317932
$ws.Range("G15").Value = 3.598404548098813
$ws.Range("I15").Value = 23.88712573867756
$ws.Range("L15").Value = 14.73990982801404
$ws.Range("N15").Value = 17.6470351910562

$ws.Range("B16").Value = 19.22572259612454
$ws.Range("D16").Value = 8.398942873665211
$ws.Range("E16").Value = 29.17157491561163
$ws.Range("F16").Value = 51.28049335547568
$ws.Range("G16").Value = 3.604096335435395
$ws.Range("I16").Value = 23.99496034121906
$ws.Range("L16").Value = 14.51428427870188
$ws.Range("N16").Value = 17.69655853685876

$ws.Range("B17").Value = 19.04769102932986
$ws.Range("D17").Value = 8.335533921725249
$ws.Range("E17").Value = 28.58865442957559
$ws.Range("F17").Value = 50.63908483037353
$ws.Range("G17").Value = 3.607651022071297
$ws.Range("I17").Value = 24.06296726834987
$ws.Range("L17").Value = 14.37541301819978
$ws.Range("N17").Value = 17.72748469211931

$ws.Range("B18").Value = 18.94497594707497
$ws.Range("D18").Value = 8.299273171675043
$ws.Range("E18").Value = 28.24822898485871
$ws.Range("F18").Value = 50.26895266529868
$ws.Range("G18").Value = 3.609718881752098
$ws.Range("I18").Value = 24.10276097208078
$ws.Range("L18").Value = 14.29539824878117
$ws.Range("N18").Value = 17.74547406143705

$ws.Range("B19").Value = 18.91014764867866
$ws.Range("D19").Value = 8.287033150516434
$ws.Range("E19").Value = 28.1320820753012
$ws.Range("F19").Value = 50.1434374621613
$ws.Range("G19").Value = 3.610423040077738
$ws.Range("I19").Value = 24.11635055629928
$ws.Range("L19").Value = 14.2682853821348
$ws.Range("N19").Value = 17.75159965896479

$ws.Range("B20").Value = 19.0666763388823
$ws.Range("D20").Value = 8.342262400714805
$ws.Range("E20").Value = 28.65123954534047
$ws.Range("F20").Value = 50.70749229291972
$ws.Range("G20").Value = 3.60727021196409
$ws.Range("I20").Value = 24.05565757803141
$ws.Range("L20").Value = 14.39021120904179
$ws.Range("N20").Value = 17.72417172439733

$ws.Range("B21").Value = 19.58856034703593
$ws.Range("D21").Value = 8.53025707888739
$ws.Range("E21").Value = 30.3341489816899
$ws.Range("F21").Value = 52.58723607831427
$ws.Range("G21").Value = 3.596954373213837
$ws.Range("I21").Value = 23.85986094019772
$ws.Range("L21").Value = 14.79799798424186
$ws.Range("N21").Value = 17.63441687953972

$ws.Range("B22").Value = 19.92607024293453
$ws.Range("D22").Value = 8.654748996525891
$ws.Range("E22").Value = 31.38847377049994
$ws.Range("F22").Value = 53.80222410980121
$ws.Range("G22").Value = 3.590416408446753
$ws.Range("I22").Value = 23.73800678495358
$ws.Range("L22").Value = 15.06267736028389
$ws.Range("N22").Value = 17.57752727106497

$ws.Range("B23").Value = 19.74629682801747
$ws.Range("D23").Value = 8.588168352496167
$ws.Range("E23").Value = 30.8299152794713
$ws.Range("F23").Value = 53.15512557652241
$ws.Range("G23").Value = 3.593887342408592
$ws.Range("I23").Value = 23.80247985525896
$ws.Range("L23").Value = 14.92160807804197
$ws.Range("N23").Value = 17.60772948760832

$ws.Range("B24").Value = 19.0580942005105
$ws.Range("D24").Value = 8.339219848649869
$ws.Range("E24").Value = 28.62296131528428
$ws.Range("F24").Value = 50.67656953923679
$ws.Range("G24").Value = 3.60744230072648
$ws.Range("I24").Value = 24.05896012313488
$ws.Range("L24").Value = 14.3835214879509
$ws.Range("N24").Value = 17.7256688642643

$ws.Range("B25").Value = 18.30419062627572
$ws.Range("D25").Value = 8.078879008332294
$ws.Range("E25").Value = 26.04559857561935
$ws.Range("F25").Value = 47.95876474912405
$ws.Range("G25").Value = 3.622956235202637
$ws.Range("I25").Value = 24.36149710526829
$ws.Range("L25").Value = 13.7981672625321
$ws.Range("N25").Value = 17.86060195689358
